# UploadTemplate.xlsx edit:
#  - "Adobe AAM" sheet: the values that used to live in J3/K3 (Modeling
#    Price / Modeling UoM columns) are moved two columns to the right,
#    into L3/M3 (Activation Price / Activation UoM), and that sheet
#    becomes the active tab/selection.
#  - "Adobe AdCloud" sheet (previously the active tab) is no longer the
#    active sheet.

$wb = $excel.ActiveWorkbook

$wsAAM = $wb.Worksheets.Item("Adobe AAM")
$wsAdCloud = $wb.Worksheets.Item("Adobe AdCloud")

# Move J3 -> L3 and K3 -> M3 on the "Adobe AAM" sheet.
$wsAAM.Range("J3").Cut($wsAAM.Range("L3"))
$wsAAM.Range("K3").Cut($wsAAM.Range("M3"))

# Make "Adobe AAM" the active sheet/tab, with L3:M3 selected
# (activeCell = L3, sqref = L3:M3) — this also clears tabSelected from
# whatever sheet was previously active ("Adobe AdCloud").
$wsAAM.Activate()
$wsAAM.Range("L3:M3").Select()
